$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compromisos")

# Give every column (A:P) a fixed 20-character width.
$ws.Range("A1:P1").ColumnWidth = 19.1640625

# The header box border goes from "medium" to "thin" (rows 1-2), and the
# new row 3 gets a matching thin box border around each of its cells.
$header1 = $ws.Range("A1:P1")
$header1.Borders.LineStyle = 1
$header1.Borders.Weight = 2

$header2 = $ws.Range("A2:P2")
$header2.Borders.LineStyle = 1
$header2.Borders.Weight = 2

$row3 = $ws.Range("A3:P3")
$row3.Borders.LineStyle = 1
$row3.Borders.Weight = 2
$ws.Range("A3").Value = 1

# Add the new "Leyenda" worksheet right after "Compromisos".
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Leyenda"

$wb.Worksheets.Item("Compromisos").Activate()
